$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Update the report date on the title page.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Date: 2024-09-23", $true, $false, $false, $false,
                         $false, $true, 1, $false, "Date: 2024-10-08", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. "Overview of the nodes in the control flow" table (1st table).
#    - Execute SQL Task occurrences: 6 -> 5
#    - Swap the "Foreach Loop Container" / "Expression Task" row labels.
# ---------------------------------------------------------------------------
$t1 = $d.Tables.Item(1)
$t1.Cell(2, 2).Range.Text = "5"
$t1.Cell(3, 1).Range.Text = "Expression Task"
$t1.Cell(4, 1).Range.Text = "Foreach Loop Container"

# ---------------------------------------------------------------------------
# 3. "Overview of the nodes in the data flow" table (4th table).
#    - Remove the "DataSources" row that currently sits right after the
#      header (value 3).
#    - "DataDestinations" occurrences: 1 -> 2
#    - Re-insert a "DataSources" row (value 2) right after the
#      "DataDestinations" row.
# ---------------------------------------------------------------------------
$t2 = $d.Tables.Item(4)
$t2.Rows.Item(2).Delete()

$t2.Cell(6, 2).Range.Text = "2"

$t2.Rows.Add($t2.Rows.Item(7)) | Out-Null
$t2.Cell(7, 1).Range.Text = "DataSources"
$t2.Cell(7, 2).Range.Text = "2"

# ---------------------------------------------------------------------------
# 4. "Overview of utilised target tables in the data flow" table (6th table).
#    - Insert a new "Supp_Prod_output" row (value 1) right after the header
#      row, before "Error_lines".
# ---------------------------------------------------------------------------
$t3 = $d.Tables.Item(6)
$t3.Rows.Add($t3.Rows.Item(2)) | Out-Null
$t3.Cell(2, 1).Range.Text = "Supp_Prod_output"
$t3.Cell(2, 2).Range.Text = "1"
